# NPJohn.pptx - "Change minor thing in presentation"
#
# Slide 1 (Title slide):
#   - Title placeholder text "Parallel JTR" -> "NPJohn"
#   - Subtitle placeholder: "Una versione multicore e distribuita di John the
#     Ripper che sfrutta MPI." -> "Una versione parallela e distribuita di
#     John the Ripper che sfrutta MPI." (only the "multicore" -> "parallela"
#     portion changes; "Ripper che sfrutta MPI." stays untouched)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape ("Title 1") ---
$title = $s.Shapes.Item("Title 1")
$title.TextFrame.TextRange.Text = "NPJohn"

# --- Subtitle shape ("Subtitle 2") ---
$subtitle = $s.Shapes.Item("Subtitle 2")
$subtitleRange = $subtitle.TextFrame.TextRange

# Replace only the leading "Una versione multicore e distribuita di John the "
# portion (49 characters) with the new wording, leaving the trailing
# "Ripper che sfrutta MPI." runs untouched.
$lead = $subtitleRange.Characters(1, 49)
$lead.Text = "Una versione parallela e distribuita di John the "
